$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    # Force text storage so numeric-looking strings (e.g. "595.78") are not
    # silently coerced into floating point numbers by the COM layer, then
    # drop back to the default "Normal" style so no stray formatting is left
    # behind on the cell.
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell $ws "D2" "67.286.23"
$ws.Range("E2").Value = "  -4.60%  "

# Row 3 - Ethereum
Set-TextCell $ws "D3" "3.262.76"
$ws.Range("E3").Value = "  -7.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextCell $ws "D5" "595.78"
$ws.Range("E5").Value = "  -4.49%  "

# Row 6 - Solana
Set-TextCell $ws "D6" "150.75"
$ws.Range("E6").Value = "  -12.36%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - LidoStakedEther
Set-TextCell $ws "D8" "3.254.65"
$ws.Range("E8").Value = "  -7.13%  "

# Row 9 - XRP
Set-TextCell $ws "D9" "0.540"
$ws.Range("E9").Value = "  -11.30%  "

# Row 10 - Dogecoin
Set-TextCell $ws "D10" "0.170"
$ws.Range("E10").Value = "  -13.92%  "

# Row 11 - Toncoin
Set-TextCell $ws "D11" "6.61"
$ws.Range("E11").Value = "  -7.74%  "

# Row 12 - Cardano
Set-TextCell $ws "D12" "0.507"
$ws.Range("E12").Value = "  -13.51%  "

# Row 13 - Avalanche
Set-TextCell $ws "D13" "37.99"
$ws.Range("E13").Value = "  -17.72%  "

# Row 14 - ShibaInu
Set-TextCell $ws "D14" "0.0000243"
$ws.Range("E14").Value = "  -11.96%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextCell $ws "D15" "3.781.90"
$ws.Range("E15").Value = "  -7.30%  "

# Row 16 - WrappedBTC
Set-TextCell $ws "D16" "67.312.04"
$ws.Range("E16").Value = "  -4.71%  "

# Row 17 - WrappedEther
Set-TextCell $ws "D17" "3.263.97"
$ws.Range("E17").Value = "  -7.10%  "

# Row 18 - was BitcoinCash, becomes TRON
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D18" "0.114"
$ws.Range("E18").Value = "  -6.30%  "

# Row 19 - was TRON, becomes BitcoinCash
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D19" "531.36"
$ws.Range("E19").Value = "  -12.38%  "

# Row 20 - Polkadot
Set-TextCell $ws "D20" "7.19"
$ws.Range("E20").Value = "  -14.21%  "

# Row 21 - Chainlink
Set-TextCell $ws "D21" "15.03"
$ws.Range("E21").Value = "  -14.98%  "

# Row 22 - Polygon
Set-TextCell $ws "D22" "0.758"
$ws.Range("E22").Value = "  -13.75%  "

# Row 23 - Uniswap
Set-TextCell $ws "D23" "7.86"
$ws.Range("E23").Value = "  -13.42%  "

# Row 24 - Litecoin
Set-TextCell $ws "D24" "85.22"
$ws.Range("E24").Value = "  -12.19%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextCell $ws "D25" "13.55"

# Row 26 - Dai
Set-TextCell $ws "D26" "0.999"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27 - PancakeSwap
Set-TextCell $ws "D27" "3.25"
$ws.Range("E27").Value = "  -12.93%  "

# Row 28 - EthereumClassic
Set-TextCell $ws "D28" "29.28"
$ws.Range("E28").Value = "  -12.56%  "

# Row 29 - RenderToken
Set-TextCell $ws "D29" "7.96"
$ws.Range("E29").Value = "  -11.54%  "

# Row 30 - ImmutableX
Set-TextCell $ws "D30" "2.12"
$ws.Range("E30").Value = "  -16.92%  "

# Row 31 - Stacks
$ws.Range("E31").Value = "  -11.78%  "

# Row 32 - Mantle
$ws.Range("E32").Value = "  -11.09%  "

# Row 33 - Bittensor
Set-TextCell $ws "D33" "541.03"
$ws.Range("E33").Value = "  -13.16%  "

# Row 34 - Filecoin
Set-TextCell $ws "D34" "6.60"
$ws.Range("E34").Value = "  -18.14%  "

# Row 35 - NEARProtocol
Set-TextCell $ws "D35" "5.67"
$ws.Range("E35").Value = "  -16.45%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.04%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -8.07%  "

# Row 38 - OKB
Set-TextCell $ws "D38" "53.29"
$ws.Range("E38").Value = "  -5.84%  "

# Row 39 - Hedera
Set-TextCell $ws "D39" "0.0856"
$ws.Range("E39").Value = "  -13.86%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -9.59%  "

# Row 41 - Cosmos
Set-TextCell $ws "D41" "9.03"
$ws.Range("E41").Value = "  -16.22%  "

# Row 42 - Maker
Set-TextCell $ws "D42" "2.929.29"
$ws.Range("E42").Value = "  -12.20%  "

# Row 43 - dogwifhat
Set-TextCell $ws "D43" "2.70"
$ws.Range("E43").Value = "  -20.79%  "

# Row 44 - TheGraph
Set-TextCell $ws "D44" "0.261"
$ws.Range("E44").Value = "  -15.93%  "

# Row 45 - PEPE (subscript-3 char, U+2083)
$sub3 = [char]8323
Set-TextCell $ws "D45" ("0.0{0}0586" -f $sub3)
$ws.Range("E45").Value = "  -18.54%  "

# Row 46 - InjectiveProtocol
Set-TextCell $ws "D46" "26.56"
$ws.Range("E46").Value = "  -16.57%  "

# Row 47 - was Fetch.AI, becomes USDe
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws "D47" "1.00"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48 - was USDe, becomes Fetch.AI
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws "D48" "2.14"
$ws.Range("E48").Value = "  -14.82%  "

# Row 49 - Monero
Set-TextCell $ws "D49" "127.43"
$ws.Range("E49").Value = "  -4.36%  "

# Row 50 - ThetaToken
Set-TextCell $ws "D50" "2.31"
$ws.Range("E50").Value = "  -21.63%  "

# Row 51 - Stellar
Set-TextCell $ws "D51" "0.113"
$ws.Range("E51").Value = "  -12.96%  "
